$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and one swapped pair of rows)

$ws.Range("D2").Value = "'36.980.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.72%  "

$ws.Range("D3").Value = "'1.988.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.45%  "

$ws.Range("D5").Value = "'239.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.00%  "

$ws.Range("D6").Value = "'0.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.70%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'54.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.58%  "

$ws.Range("D9").Value = "'0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.18%  "

$ws.Range("D10").Value = "'58.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Value = "'0.0750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.29%  "

$ws.Range("D12").Value = "'0.0982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.29%  "

$ws.Range("D13").Value = "'2.285.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.36%  "

$ws.Range("D14").Value = "'14.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.41%  "

$ws.Range("D15").Value = "'20.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").Value = "'0.754"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.18%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.026.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'5.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.20%  "

$ws.Range("D19").Value = "'36.960.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").Value = "'68.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").Value = "'0.0₃0806"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.94%  "

$ws.Range("D22").Value = "'228.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").Value = "'4.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.12%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").Value = "'2.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.77%  "

$ws.Range("D26").Value = "'2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "'161.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").Value = "'8.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.68%  "

$ws.Range("D29").Value = "'19.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.04%  "

$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.42%  "

$ws.Range("D31").Value = "'1.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.76%  "

$ws.Range("D32").Value = "'0.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.89%  "

$ws.Range("D33").Value = "'4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.54%  "

$ws.Range("D34").Value = "'0.0612"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.85%  "

$ws.Range("D35").Value = "'4.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.09%  "

$ws.Range("D36").Value = "'2.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.39%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").Value = "'3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.07%  "

$ws.Range("D40").Value = "'5.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.40%  "

$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").Value = "'1.429.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("E43").Value = "  -6.57%  "

$ws.Range("D44").Value = "'1.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.92%  "

$ws.Range("D45").Value = "'0.0884"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.76%  "

$ws.Range("D46").Value = "'88.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.22%  "

$ws.Range("D47").Value = "'15.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.45%  "

$ws.Range("D48").Value = "'0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.06%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  +13.50%  "

$ws.Range("D51").Value = "'6.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.51%  "
